$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 5).Value = 3
$ws.Cells.Item(2, 6).Value = 1
$ws.Cells.Item(2, 7).Value = 1.824961333333333
$ws.Cells.Item(2, 8).Value = 5.474884
$ws.Cells.Item(2, 9).Value = 0.377840167393297
$ws.Cells.Item(2, 10).Value = 0.3778401673932969
$ws.Cells.Item(2, 13).Value = 0.7489546666666667
$ws.Cells.Item(2, 14).Value = 2.246864
$ws.Cells.Item(2, 15).Value = 0.05220789806691288
$ws.Cells.Item(2, 16).Value = 0.05220789806691287
$ws.Cells.Item(2, 17).Value = 1.366813307086222
$ws.Cells.Item(2, 18).Value = 12.301319763776
$ws.Cells.Item(2, 19).Value = 0.01972624094485455
$ws.Cells.Item(2, 20).Value = 0.01972624094485454
$ws.Cells.Item(3, 5).Value = 3
$ws.Cells.Item(3, 6).Value = 1
$ws.Cells.Item(3, 7).Value = 1.824961333333333
$ws.Cells.Item(3, 8).Value = 5.474884
$ws.Cells.Item(3, 9).Value = 0.377840167393297
$ws.Cells.Item(3, 10).Value = 0.3778401673932969
$ws.Cells.Item(3, 15).Value = 0.1982273102638064
$ws.Cells.Item(3, 16).Value = 0.1982273102638064
$ws.Cells.Item(3, 17).Value = 5.189630985511556
$ws.Cells.Item(3, 18).Value = 46.706678869604
$ws.Cells.Item(3, 19).Value = 0.07489824009199963
$ws.Cells.Item(3, 20).Value = 0.07489824009199962
$ws.Cells.Item(4, 5).Value = 3
$ws.Cells.Item(4, 6).Value = 1
$ws.Cells.Item(4, 7).Value = 1.824961333333333
$ws.Cells.Item(4, 8).Value = 5.474884
$ws.Cells.Item(4, 9).Value = 0.377840167393297
$ws.Cells.Item(4, 10).Value = 0.3778401673932969
$ws.Cells.Item(4, 13).Value = 10.337765
$ws.Cells.Item(4, 14).Value = 31.013295
$ws.Cells.Item(4, 15).Value = 0.7206216949842531
$ws.Cells.Item(4, 16).Value = 0.720621694984253
$ws.Cells.Item(4, 17).Value = 18.86602139808667
$ws.Cells.Item(4, 18).Value = 169.79419258278
$ws.Cells.Item(4, 19).Value = 0.2722798218600916
$ws.Cells.Item(4, 20).Value = 0.2722798218600915
$ws.Cells.Item(5, 5).Value = 3
$ws.Cells.Item(5, 6).Value = 1
$ws.Cells.Item(5, 7).Value = 1.824961333333333
$ws.Cells.Item(5, 8).Value = 5.474884
$ws.Cells.Item(5, 9).Value = 0.377840167393297
$ws.Cells.Item(5, 10).Value = 0.3778401673932969
$ws.Cells.Item(5, 13).Value = 0.4152066666666667
$ws.Cells.Item(5, 14).Value = 1.24562
$ws.Cells.Item(5, 15).Value = 0.02894309668502767
$ws.Cells.Item(5, 16).Value = 0.02894309668502767
$ws.Cells.Item(5, 17).Value = 0.7577361120088889
$ws.Cells.Item(5, 18).Value = 6.81962500808
$ws.Cells.Item(5, 19).Value = 0.01093586449635123
$ws.Cells.Item(5, 20).Value = 0.01093586449635123
$ws.Cells.Item(6, 9).Value = 0.03077064395059555
$ws.Cells.Item(6, 10).Value = 0.03077064395059554
$ws.Cells.Item(6, 13).Value = 0.7489546666666667
$ws.Cells.Item(6, 14).Value = 2.246864
$ws.Cells.Item(6, 15).Value = 0.05220789806691288
$ws.Cells.Item(6, 16).Value = 0.05220789806691287
$ws.Cells.Item(6, 17).Value = 0.1113108908177778
$ws.Cells.Item(6, 18).Value = 1.00179801736
$ws.Cells.Item(6, 19).Value = 0.001606470642825962
$ws.Cells.Item(6, 20).Value = 0.001606470642825961
$ws.Cells.Item(7, 9).Value = 0.03077064395059555
$ws.Cells.Item(7, 10).Value = 0.03077064395059554
$ws.Cells.Item(7, 15).Value = 0.1982273102638064
$ws.Cells.Item(7, 16).Value = 0.1982273102638064
$ws.Cells.Item(7, 19).Value = 0.006099581985411822
$ws.Cells.Item(7, 20).Value = 0.00609958198541182
$ws.Cells.Item(8, 9).Value = 0.03077064395059555
$ws.Cells.Item(8, 10).Value = 0.03077064395059554
$ws.Cells.Item(8, 13).Value = 10.337765
$ws.Cells.Item(8, 14).Value = 31.013295
$ws.Cells.Item(8, 15).Value = 0.7206216949842531
$ws.Cells.Item(8, 16).Value = 0.720621694984253
$ws.Cells.Item(8, 17).Value = 1.536415863908333
$ws.Cells.Item(8, 18).Value = 13.827742775175
$ws.Cells.Item(8, 19).Value = 0.02217399359943512
$ws.Cells.Item(8, 20).Value = 0.02217399359943511
$ws.Cells.Item(9, 9).Value = 0.03077064395059555
$ws.Cells.Item(9, 10).Value = 0.03077064395059554
$ws.Cells.Item(9, 13).Value = 0.4152066666666667
$ws.Cells.Item(9, 14).Value = 1.24562
$ws.Cells.Item(9, 15).Value = 0.02894309668502767
$ws.Cells.Item(9, 16).Value = 0.02894309668502767
$ws.Cells.Item(9, 17).Value = 0.06170870681111112
$ws.Cells.Item(9, 18).Value = 0.5553783613
$ws.Cells.Item(9, 19).Value = 0.0008905977229226489
$ws.Cells.Item(9, 20).Value = 0.0008905977229226484
$ws.Cells.Item(10, 7).Value = 2.658767
$ws.Cells.Item(10, 8).Value = 7.976300999999999
$ws.Cells.Item(10, 9).Value = 0.5504713716344166
$ws.Cells.Item(10, 10).Value = 0.5504713716344165
$ws.Cells.Item(10, 13).Value = 0.7489546666666667
$ws.Cells.Item(10, 14).Value = 2.246864
$ws.Cells.Item(10, 15).Value = 0.05220789806691288
$ws.Cells.Item(10, 16).Value = 0.05220789806691287
$ws.Cells.Item(10, 17).Value = 1.991295952229333
$ws.Cells.Item(10, 18).Value = 17.921663570064
$ws.Cells.Item(10, 19).Value = 0.02873895325904334
$ws.Cells.Item(10, 20).Value = 0.02873895325904333
$ws.Cells.Item(11, 7).Value = 2.658767
$ws.Cells.Item(11, 8).Value = 7.976300999999999
$ws.Cells.Item(11, 9).Value = 0.5504713716344166
$ws.Cells.Item(11, 10).Value = 0.5504713716344165
$ws.Cells.Item(11, 15).Value = 0.1982273102638064
$ws.Cells.Item(11, 16).Value = 0.1982273102638064
$ws.Cells.Item(11, 17).Value = 7.560718879042333
$ws.Cells.Item(11, 18).Value = 68.046469911381
$ws.Cells.Item(11, 19).Value = 0.1091184593763186
$ws.Cells.Item(11, 20).Value = 0.1091184593763186
$ws.Cells.Item(12, 7).Value = 2.658767
$ws.Cells.Item(12, 8).Value = 7.976300999999999
$ws.Cells.Item(12, 9).Value = 0.5504713716344166
$ws.Cells.Item(12, 10).Value = 0.5504713716344165
$ws.Cells.Item(12, 13).Value = 10.337765
$ws.Cells.Item(12, 14).Value = 31.013295
$ws.Cells.Item(12, 15).Value = 0.7206216949842531
$ws.Cells.Item(12, 16).Value = 0.720621694984253
$ws.Cells.Item(12, 17).Value = 27.48570843575499
$ws.Cells.Item(12, 18).Value = 247.371375921795
$ws.Cells.Item(12, 19).Value = 0.3966816128675
$ws.Cells.Item(12, 20).Value = 0.3966816128674999
$ws.Cells.Item(13, 7).Value = 2.658767
$ws.Cells.Item(13, 8).Value = 7.976300999999999
$ws.Cells.Item(13, 9).Value = 0.5504713716344166
$ws.Cells.Item(13, 10).Value = 0.5504713716344165
$ws.Cells.Item(13, 13).Value = 0.4152066666666667
$ws.Cells.Item(13, 14).Value = 1.24562
$ws.Cells.Item(13, 15).Value = 0.02894309668502767
$ws.Cells.Item(13, 16).Value = 0.02894309668502767
$ws.Cells.Item(13, 17).Value = 1.103937783513333
$ws.Cells.Item(13, 18).Value = 9.935440051619999
$ws.Cells.Item(13, 19).Value = 0.01593234613155472
$ws.Cells.Item(13, 20).Value = 0.01593234613155472
$ws.Cells.Item(14, 5).Value = 2
$ws.Cells.Item(14, 6).Value = 0.6666666666666666
$ws.Cells.Item(14, 7).Value = 0.1976323333333333
$ws.Cells.Item(14, 8).Value = 0.592897
$ws.Cells.Item(14, 9).Value = 0.04091781702169097
$ws.Cells.Item(14, 10).Value = 0.04091781702169097
$ws.Cells.Item(14, 13).Value = 0.7489546666666667
$ws.Cells.Item(14, 14).Value = 2.246864
$ws.Cells.Item(14, 15).Value = 0.05220789806691288
$ws.Cells.Item(14, 16).Value = 0.05220789806691287
$ws.Cells.Item(14, 17).Value = 0.1480176583342222
$ws.Cells.Item(14, 18).Value = 1.332158925008
$ws.Cells.Item(14, 19).Value = 0.002136233220189035
$ws.Cells.Item(14, 20).Value = 0.002136233220189035
$ws.Cells.Item(15, 5).Value = 2
$ws.Cells.Item(15, 6).Value = 0.6666666666666666
$ws.Cells.Item(15, 7).Value = 0.1976323333333333
$ws.Cells.Item(15, 8).Value = 0.592897
$ws.Cells.Item(15, 9).Value = 0.04091781702169097
$ws.Cells.Item(15, 10).Value = 0.04091781702169097
$ws.Cells.Item(15, 15).Value = 0.1982273102638064
$ws.Cells.Item(15, 16).Value = 0.1982273102638064
$ws.Cells.Item(15, 17).Value = 0.5620058146285556
$ws.Cells.Item(15, 18).Value = 5.058052331657001
$ws.Cells.Item(15, 19).Value = 0.008111028810076396
$ws.Cells.Item(15, 20).Value = 0.008111028810076396
$ws.Cells.Item(16, 5).Value = 2
$ws.Cells.Item(16, 6).Value = 0.6666666666666666
$ws.Cells.Item(16, 7).Value = 0.1976323333333333
$ws.Cells.Item(16, 8).Value = 0.592897
$ws.Cells.Item(16, 9).Value = 0.04091781702169097
$ws.Cells.Item(16, 10).Value = 0.04091781702169097
$ws.Cells.Item(16, 13).Value = 10.337765
$ws.Cells.Item(16, 14).Value = 31.013295
$ws.Cells.Item(16, 15).Value = 0.7206216949842531
$ws.Cells.Item(16, 16).Value = 0.720621694984253
$ws.Cells.Item(16, 17).Value = 2.043076618401666
$ws.Cells.Item(16, 18).Value = 18.387689565615
$ws.Cells.Item(16, 19).Value = 0.02948626665722647
$ws.Cells.Item(16, 20).Value = 0.02948626665722646
$ws.Cells.Item(17, 5).Value = 2
$ws.Cells.Item(17, 6).Value = 0.6666666666666666
$ws.Cells.Item(17, 7).Value = 0.1976323333333333
$ws.Cells.Item(17, 8).Value = 0.592897
$ws.Cells.Item(17, 9).Value = 0.04091781702169097
$ws.Cells.Item(17, 10).Value = 0.04091781702169097
$ws.Cells.Item(17, 13).Value = 0.4152066666666667
$ws.Cells.Item(17, 14).Value = 1.24562
$ws.Cells.Item(17, 15).Value = 0.02894309668502767
$ws.Cells.Item(17, 16).Value = 0.02894309668502767
$ws.Cells.Item(17, 17).Value = 0.08205826234888888
$ws.Cells.Item(17, 18).Value = 0.73852436114
$ws.Cells.Item(17, 19).Value = 0.001184288334199073
$ws.Cells.Item(17, 20).Value = 0.001184288334199073
